# Update the "Last Updated: 20-Apr-20" date stamp to "21-Apr-20" on every
# slide. The date text lives in its own run at the end of the first shape's
# (TextBox 2) text body on each slide, so we locate the "20-Apr-20"
# substring within the shape's full text and replace just that slice,
# leaving every other run (and its formatting) untouched.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $shape = $slide.Shapes.Item(1)

    if (-not $shape.HasTextFrame) { continue }

    $textRange = $shape.TextFrame.TextRange
    $fullText = $textRange.Text
    $searchTerm = "20-Apr-20"
    $pos = $fullText.IndexOf($searchTerm)

    if ($pos -ge 0) {
        $target = $textRange.Characters($pos + 1, $searchTerm.Length)
        $target.Text = "21-Apr-20"
    }
}
